# Sync attendance_reports: reorder "Recorded By" entries so that "System"
# appears at the end of the list instead of the beginning.
#
# For the "Recorded By" column (G) on the active sheet, several values
# start with "System, ..." - the word "System" needs to move from the
# front of the comma-separated list to the back.
#
# Example:
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com"         -> "backup@backdoor.com, System"
#   "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row/column on the sheet.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Locate the "Recorded By" column by scanning the header row (row 1).
$lastCol = $usedRange.Columns.Count + $usedRange.Column - 1
$recordedByCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $headerValue = $ws.Cells.Item(1, $c).Value()
    if ($headerValue -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}

if ($recordedByCol -eq 0) {
    $recordedByCol = 7  # fall back to column G
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $value = $cell.Value()

    if ($value -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($value -eq "System, backup@backdoor.com") {
        $cell.Value = "backup@backdoor.com, System"
    }
    elseif ($value -eq "System, system, backup@backdoor.com") {
        $cell.Value = "system, backup@backdoor.com, System"
    }
}
